# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets
# to reflect refreshed source data at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$f1 = @{
    2 = 36
    3 = 235
    4 = 4841
    6 = 162
    7 = 124
    9 = 97
    10 = 767
    11 = 230
    12 = 1204
    13 = 116
    14 = 258
    15 = 193
    18 = 155
    20 = 4082
    21 = 6390
    23 = 0
    25 = 544
    27 = 3995
    29 = 49
    30 = 26
    31 = 2595
    32 = 569
    34 = 147
    35 = 303
    36 = 316
    38 = 183
    42 = 48
    43 = 77
    44 = 60
    45 = 0
    47 = 3
    48 = 79
}
foreach ($row in $f1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $f1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
$f4 = @{
    3 = 0
    4 = 4841
    6 = 162
    7 = 0
    8 = 112
    10 = 97
    11 = 767
    12 = 230
    13 = 1204
    14 = 116
    15 = 193
    16 = 0
    18 = 155
    20 = 4082
    21 = 6390
    23 = 40
    25 = 544
    26 = 48
    27 = 3995
    29 = 49
    30 = 26
    31 = 2595
    32 = 0
    34 = 147
    35 = 303
    36 = 316
    38 = 183
    39 = 11
    41 = 977
    42 = 48
    43 = 77
    45 = 502
    47 = 3
    48 = 79
    49 = 590
}
foreach ($row in $f4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $f4[$row]
}
